$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.480.48'
$ws.Range('E2').Value = '  +0.66%  '
$ws.Range('D3').Value = '1.897.43'
$ws.Range('E3').Value = '  +0.58%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('B5').Value = 'XRP'
$ws.Range('C5').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.693'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.74%  '
$ws.Range('B6').Value = 'BNB'
$ws.Range('C6').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '246.46'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.15%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '43.12'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -1.59%  '
$ws.Range('E9').Value = '  +2.21%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '56.23'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +8.81%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0757'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +2.73%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0983'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +1.41%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.23'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +8.75%  '
$ws.Range('E14').Value = '  +10.15%  '
$ws.Range('D15').Value = '2.172.06'
$ws.Range('E15').Value = '  +0.44%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.02'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +2.34%  '
$ws.Range('D17').Value = '1.898.98'
$ws.Range('E17').Value = '  +0.36%  '
$ws.Range('D18').Value = '35.443.18'
$ws.Range('E18').Value = '  +0.62%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '73.71'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.98%  '
$ws.Range('D20').Value = '0.0₃0832'
$ws.Range('E20').Value = '  +1.59%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '244.57'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.28%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '13.01'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +1.58%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.21'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +5.03%  '
$ws.Range('E24').Value = '  +5.51%  '
$ws.Range('E25').Value = '  -0.06%  '
$ws.Range('E26').Value = '  -0.94%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '166.69'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +1.02%  '
$ws.Range('E28').Value = '  +1.77%  '
$ws.Range('E29').Value = '  +0.17%  '
$ws.Range('E30').Value = '  +0.81%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.37'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +2.65%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0603'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +4.21%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.25'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.60%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.87'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +21.85%  '
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.48'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -16.06%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.855'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.57%  '
$ws.Range('B38').Value = 'LidoDAOToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.95'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -1.75%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0736'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +10.14%  '
$ws.Range('E40').Value = '  +6.10%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '99.04'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +1.48%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '16.98'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -1.17%  '
$ws.Range('E43').Value = '  -0.82%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.75'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +14.59%  '
$ws.Range('D45').Value = '1.325.40'
$ws.Range('E45').Value = '  +2.54%  '
$ws.Range('E46').Value = '  +0.94%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0812'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.16%  '
$ws.Range('E48').Value = '  +0.74%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.74'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.17%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.40'
$ws.Range('D50').ClearFormats()
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '42.56'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.11%  '
